# "ajuste de cronograma" - burndown sprint update
# Fill in the burndown diagonal marker cells for the next few days and
# switch the active sheet over to the "Gráfico" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 1")

$ws.Range("D4").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("G7").Value = 1

$ws.Range("H8").Select()

$chartSheet = $wb.Worksheets.Item("Gráfico")
$chartSheet.Activate()
